# Natmi LR-pairs re-run ("Natmi following Dr Hou advice"): a 3rd cluster
# ("ECs") was added to the Young-D0 cluster set, so the Wnt5a-Ror1 sending/
# target-cluster cross product grows from 2x2 to 2x3 rows, and every row's
# computed NATMI statistics (detection rate, specificity, edge weights, ...)
# is refreshed with the recomputed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Ror1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9.156959333333335
$ws.Range("H2").Value = 27.470878
$ws.Range("I2").Value = 0.969469463764299
$ws.Range("J2").Value = 0.9694694637642989
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.921193
$ws.Range("N2").Value = 2.763579
$ws.Range("O2").Value = 0.04194147971775762
$ws.Range("P2").Value = 0.04194147971775761
$ws.Range("Q2").Value = 8.435326839151335
$ws.Range("R2").Value = 75.91794155236201
$ws.Range("S2").Value = 0.0406609838514557
$ws.Range("T2").Value = 0.04066098385145569

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Ror1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 9.156959333333335
$ws.Range("H3").Value = 27.470878
$ws.Range("I3").Value = 0.969469463764299
$ws.Range("J3").Value = 0.9694694637642989
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 9.502892666666668
$ws.Range("N3").Value = 28.508678
$ws.Range("O3").Value = 0.4326621891818844
$ws.Range("P3").Value = 0.4326621891818844
$ws.Range("Q3").Value = 87.01760169769824
$ws.Range("R3").Value = 783.1584152792842
$ws.Range("S3").Value = 0.4194527805372491
$ws.Range("T3").Value = 0.4194527805372491

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Ror1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 9.156959333333335
$ws.Range("H4").Value = 27.470878
$ws.Range("I4").Value = 0.969469463764299
$ws.Range("J4").Value = 0.9694694637642989
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 11.539684
$ws.Range("N4").Value = 34.619052
$ws.Range("O4").Value = 0.525396331100358
$ws.Range("P4").Value = 0.5253963311003579
$ws.Range("Q4").Value = 105.6684171075173
$ws.Range("R4").Value = 951.015753967656
$ws.Range("S4").Value = 0.5093556993755941
$ws.Range("T4").Value = 0.509355699375594

$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Ror1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.288371
$ws.Range("H5").Value = 0.865113
$ws.Range("I5").Value = 0.03053053623570109
$ws.Range("J5").Value = 0.03053053623570109
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.921193
$ws.Range("N5").Value = 2.763579
$ws.Range("O5").Value = 0.04194147971775762
$ws.Range("P5").Value = 0.04194147971775761
$ws.Range("Q5").Value = 0.265645346603
$ws.Range("R5").Value = 2.390808119427
$ws.Range("S5").Value = 0.001280495866301921
$ws.Range("T5").Value = 0.001280495866301921

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Ror1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.288371
$ws.Range("H6").Value = 0.865113
$ws.Range("I6").Value = 0.03053053623570109
$ws.Range("J6").Value = 0.03053053623570109
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 9.502892666666668
$ws.Range("N6").Value = 28.508678
$ws.Range("O6").Value = 0.4326621891818844
$ws.Range("P6").Value = 0.4326621891818844
$ws.Range("Q6").Value = 2.740358661179334
$ws.Range("R6").Value = 24.663227950614
$ws.Range("S6").Value = 0.01320940864463528
$ws.Range("T6").Value = 0.01320940864463528

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Ror1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.288371
$ws.Range("H7").Value = 0.865113
$ws.Range("I7").Value = 0.03053053623570109
$ws.Range("J7").Value = 0.03053053623570109
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 11.539684
$ws.Range("N7").Value = 34.619052
$ws.Range("O7").Value = 0.525396331100358
$ws.Range("P7").Value = 0.5253963311003579
$ws.Range("Q7").Value = 3.327710214764
$ws.Range("R7").Value = 29.949391932876
$ws.Range("S7").Value = 0.01604063172476389
$ws.Range("T7").Value = 0.01604063172476388
